{"js": "// Update division problems in the table (3-digit \u00f7 1-digit) to match the\n// new worksheet data set. Each table cell's text is replaced positionally\n// (row, col) rather than by global find/replace, because some new values\n// coincide with other cells' old values (e.g. \"913\u00f77=\" is both an old value\n// and a new value at a different cell) \u2014 positional replacement avoids any\n// cross-talk between edits.\nconst entries = [\n  {\n    \"row\": 0,\n    \"col\": 0,\n    \"old\": \"608\u00f72=\",\n    \"new\": \"338\u00f72=\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 1,\n    \"old\": \"849\u00f76=\",\n    \"new\": \"549\u00f73=\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 2,\n    \"old\": \"116\u00f79=\",\n    \"new\": \"768\u00f73=\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 3,\n    \"old\": \"410\u00f77=\",\n    \"new\": \"565\u00f74=\"\n  },\n  {\n    \"row\": 0,\n    \"col\": 4,\n    \"old\": \"489\u00f75=\",\n    \"new\": \"211\u00f79=\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 0,\n    \"old\": \"545\u00f79=\",\n    \"new\": \"309\u00f74=\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 1,\n    \"old\": \"459\u00f76=\",\n    \"new\": \"850\u00f77=\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 2,\n    \"old\": \"867\u00f75=\",\n    \"new\": \"801\u00f79=\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 3,\n    \"old\": \"637\u00f75=\",\n    \"new\": \"352\u00f78=\"\n  },\n  {\n    \"row\": 4,\n    \"col\": 4,\n    \"old\": \"519\u00f73=\",\n    \"new\": \"968\u00f74=\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 0,\n    \"old\": \"318\u00f74=\",\n    \"new\": \"492\u00f73=\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 1,\n    \"old\": \"551\u00f73=\",\n    \"new\": \"412\u00f76=\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 2,\n    \"old\": \"261\u00f78=\",\n    \"new\": \"230\u00f73=\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 3,\n    \"old\": \"788\u00f74=\",\n    \"new\": \"913\u00f77=\"\n  },\n  {\n    \"row\": 8,\n    \"col\": 4,\n    \"old\": \"623\u00f78=\",\n    \"new\": \"796\u00f74=\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 0,\n    \"old\": \"690\u00f78=\",\n    \"new\": \"811\u00f74=\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 1,\n    \"old\": \"312\u00f72=\",\n    \"new\": \"320\u00f76=\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 2,\n    \"old\": \"618\u00f77=\",\n    \"new\": \"493\u00f76=\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 3,\n    \"old\": \"582\u00f73=\",\n    \"new\": \"514\u00f76=\"\n  },\n  {\n    \"row\": 12,\n    \"col\": 4,\n    \"old\": \"283\u00f75=\",\n    \"new\": \"658\u00f78=\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 0,\n    \"old\": \"859\u00f72=\",\n    \"new\": \"691\u00f76=\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 1,\n    \"old\": \"130\u00f73=\",\n    \"new\": \"675\u00f76=\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 2,\n    \"old\": \"879\u00f77=\",\n    \"new\": \"234\u00f76=\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 3,\n    \"old\": \"913\u00f77=\",\n    \"new\": \"860\u00f73=\"\n  },\n  {\n    \"row\": 16,\n    \"col\": 4,\n    \"old\": \"518\u00f78=\",\n    \"new\": \"194\u00f75=\"\n  }\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected at least one table in the document.\");\n}\n\nconst table = tables.items[0];\n\n// Load the first paragraph of each target cell so we can verify + replace\n// its text while preserving the existing run formatting (font/size).\nconst paragraphs = entries.map(e => {\n  const cell = table.getCell(e.row, e.col);\n  const para = cell.body.paragraphs.getFirst();\n  para.load(\"text\");\n  return para;\n});\nawait context.sync();\n\nentries.forEach((e, i) => {\n  const para = paragraphs[i];\n  if (para.text !== e.old) {\n    throw new Error(\n      `Unexpected text at row ${e.row}, col ${e.col}: expected \"${e.old}\" but found \"${para.text}\"`\n    );\n  }\n  para.insertText(e.new, Word.InsertLocation.replace);\n});\n\nawait context.sync();\n", "ps1": "# Update division problems in the worksheet table (3-digit / 1-digit).\n# Each table cell is addressed by its 1-based (row, col) position and its\n# text is replaced directly. Positional addressing (rather than a global\n# Find/Replace across the whole document) is used deliberately: some of the\n# new values coincide with other cells' OLD values (e.g. \"913\u00f77=\" is both\n# an old value at one cell and the new value of a different, earlier cell),\n# so a document-wide textual replace could cross-talk between edits if not\n# scoped per cell.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$cell = $t.Cell(1, 1)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"608\u00f72=\") {\n    throw \"Unexpected text in cell (1,1): expected `\"608\u00f72=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"338\u00f72=\"\n\n$cell = $t.Cell(1, 2)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"849\u00f76=\") {\n    throw \"Unexpected text in cell (1,2): expected `\"849\u00f76=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"549\u00f73=\"\n\n$cell = $t.Cell(1, 3)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"116\u00f79=\") {\n    throw \"Unexpected text in cell (1,3): expected `\"116\u00f79=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"768\u00f73=\"\n\n$cell = $t.Cell(1, 4)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"410\u00f77=\") {\n    throw \"Unexpected text in cell (1,4): expected `\"410\u00f77=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"565\u00f74=\"\n\n$cell = $t.Cell(1, 5)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"489\u00f75=\") {\n    throw \"Unexpected text in cell (1,5): expected `\"489\u00f75=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"211\u00f79=\"\n\n$cell = $t.Cell(5, 1)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"545\u00f79=\") {\n    throw \"Unexpected text in cell (5,1): expected `\"545\u00f79=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"309\u00f74=\"\n\n$cell = $t.Cell(5, 2)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"459\u00f76=\") {\n    throw \"Unexpected text in cell (5,2): expected `\"459\u00f76=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"850\u00f77=\"\n\n$cell = $t.Cell(5, 3)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"867\u00f75=\") {\n    throw \"Unexpected text in cell (5,3): expected `\"867\u00f75=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"801\u00f79=\"\n\n$cell = $t.Cell(5, 4)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"637\u00f75=\") {\n    throw \"Unexpected text in cell (5,4): expected `\"637\u00f75=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"352\u00f78=\"\n\n$cell = $t.Cell(5, 5)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"519\u00f73=\") {\n    throw \"Unexpected text in cell (5,5): expected `\"519\u00f73=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"968\u00f74=\"\n\n$cell = $t.Cell(9, 1)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"318\u00f74=\") {\n    throw \"Unexpected text in cell (9,1): expected `\"318\u00f74=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"492\u00f73=\"\n\n$cell = $t.Cell(9, 2)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"551\u00f73=\") {\n    throw \"Unexpected text in cell (9,2): expected `\"551\u00f73=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"412\u00f76=\"\n\n$cell = $t.Cell(9, 3)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"261\u00f78=\") {\n    throw \"Unexpected text in cell (9,3): expected `\"261\u00f78=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"230\u00f73=\"\n\n$cell = $t.Cell(9, 4)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"788\u00f74=\") {\n    throw \"Unexpected text in cell (9,4): expected `\"788\u00f74=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"913\u00f77=\"\n\n$cell = $t.Cell(9, 5)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"623\u00f78=\") {\n    throw \"Unexpected text in cell (9,5): expected `\"623\u00f78=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"796\u00f74=\"\n\n$cell = $t.Cell(13, 1)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"690\u00f78=\") {\n    throw \"Unexpected text in cell (13,1): expected `\"690\u00f78=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"811\u00f74=\"\n\n$cell = $t.Cell(13, 2)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"312\u00f72=\") {\n    throw \"Unexpected text in cell (13,2): expected `\"312\u00f72=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"320\u00f76=\"\n\n$cell = $t.Cell(13, 3)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"618\u00f77=\") {\n    throw \"Unexpected text in cell (13,3): expected `\"618\u00f77=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"493\u00f76=\"\n\n$cell = $t.Cell(13, 4)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"582\u00f73=\") {\n    throw \"Unexpected text in cell (13,4): expected `\"582\u00f73=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"514\u00f76=\"\n\n$cell = $t.Cell(13, 5)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"283\u00f75=\") {\n    throw \"Unexpected text in cell (13,5): expected `\"283\u00f75=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"658\u00f78=\"\n\n$cell = $t.Cell(17, 1)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"859\u00f72=\") {\n    throw \"Unexpected text in cell (17,1): expected `\"859\u00f72=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"691\u00f76=\"\n\n$cell = $t.Cell(17, 2)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"130\u00f73=\") {\n    throw \"Unexpected text in cell (17,2): expected `\"130\u00f73=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"675\u00f76=\"\n\n$cell = $t.Cell(17, 3)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"879\u00f77=\") {\n    throw \"Unexpected text in cell (17,3): expected `\"879\u00f77=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"234\u00f76=\"\n\n$cell = $t.Cell(17, 4)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"913\u00f77=\") {\n    throw \"Unexpected text in cell (17,4): expected `\"913\u00f77=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"860\u00f73=\"\n\n$cell = $t.Cell(17, 5)\n$current = $cell.Range.Text\n$currentTrimmed = $current.TrimEnd([char]13, [char]7)\nif ($currentTrimmed -ne \"518\u00f78=\") {\n    throw \"Unexpected text in cell (17,5): expected `\"518\u00f78=`\" but found `\"$currentTrimmed`\"\"\n}\n$cell.Range.Text = \"194\u00f75=\"\n\n"}
